# Update 06 Jully : add new uploads
# Adds 4 new data rows (3-6) to the "Data" sheet, each holding a
# deceased-person name / date string pair, reusing shared strings where
# the text repeats (rows 3 & 5 share the same pair; B4 & B6 share text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "સ્વ. વશરામભાઇ ભુરાભાઈ દેસાઈ"
$ws.Range("B3").Value = "સ્વર્ગવાસ તારીખ : 11/06/2025 ને મંગળવાર"

$ws.Range("A4").Value = "સ્વ. વશરામભાઇ ભુરાભાઈ દેસાઈ સ્વર્ગવાસ"
$ws.Range("B4").Value = "તારીખ : 11/06/2025 ને બુધવાર"

$ws.Range("A5").Value = "સ્વ. વશરામભાઇ ભુરાભાઈ દેસાઈ"
$ws.Range("B5").Value = "સ્વર્ગવાસ તારીખ : 11/06/2025 ને મંગળવાર"

$ws.Range("A6").Value = "સ્વ. વશરામભાઇ ભુરાભાઈ દેસાઈ "
$ws.Range("B6").Value = "તારીખ : 11/06/2025 ને બુધવાર"
